# Scene 74 — reveal the mystery speaker's name.
#
# The script at "???: Pro...!!" originally hid the speaker's identity
# behind a placeholder. This edit fills in the name ("Mara") while
# keeping the leading "?" that was already there, turning:
#     ???: Pro…!!
# into:
#     ?Mara: Pro…!!

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "???: Pro",    # FindText
    $true,          # MatchCase
    $false,         # MatchWholeWord
    $false,         # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    "?Mara: Pro",  # ReplaceWith
    2               # Replace (wdReplaceAll)
)

Write-Output "Replaced '???: Pro' -> '?Mara: Pro': $found"
